$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'" + '26.272.78'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.55%  '
# Row 3
$ws.Range('D3').Value = "'" + '1.591.59'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.22%  '
# Row 4
$ws.Range('E4').Value = '  -0.09%  '
# Row 5
$ws.Range('D5').Value = "'" + '212.57'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.82%  '
# Row 6
$ws.Range('D6').Value = "'" + '0.502'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.74%  '
# Row 7
$ws.Range('E7').Value = '  -0.10%  '
# Row 8
$ws.Range('E8').Value = '  +0.67%  '
# Row 9
$ws.Range('E9').Value = '  +0.05%  '
# Row 10
$ws.Range('D10').Value = "'" + '19.38'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.17%  '
# Row 11
$ws.Range('D11').Value = "'" + '0.0851'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.68%  '
# Row 12
$ws.Range('D12').Value = "'" + '1.815.94'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.37%  '
# Row 13
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'" + '4.02'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.70%  '
# Row 14
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = "'" + '1.568.93'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.10%  '
# Row 15
$ws.Range('E15').Value = '  +1.16%  '
# Row 16
$ws.Range('D16').Value = "'" + '64.32'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.20%  '
# Row 17
$ws.Range('D17').Value = "'" + '26.274.25'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.64%  '
# Row 18
$ws.Range('E18').Value = '  +0.57%  '
# Row 19
$ws.Range('E19').Value = '  +1.49%  '
# Row 20
$ws.Range('D20').Value = "'" + '213.97'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.01%  '
# Row 21
$ws.Range('E21').Value = '  -0.15%  '
# Row 22
$ws.Range('D22').Value = "'" + '4.30'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.30%  '
# Row 23
$ws.Range('D23').Value = "'" + '9.05'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.66%  '
# Row 24
$ws.Range('D24').Value = "'" + '2.17'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.85%  '
# Row 25
$ws.Range('D25').Value = "'" + '143.74'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.02%  '
# Row 26
$ws.Range('E26').Value = '  -0.12%  '
# Row 27
$ws.Range('D27').Value = "'" + '7.07'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.36%  '
# Row 28
$ws.Range('D28').Value = "'" + '0.112'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.12%  '
# Row 29
$ws.Range('D29').Value = "'" + '15.20'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.02%  '
# Row 30
$ws.Range('D30').Value = "'" + '0.0498'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.51%  '
# Row 31
$ws.Range('E31').Value = '  +1.80%  '
# Row 32
$ws.Range('E32').Value = '  -0.08%  '
# Row 33
$ws.Range('D33').Value = "'" + '1.347.14'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.44%  '
# Row 34
$ws.Range('E34').Value = '  -1.16%  '
# Row 35
$ws.Range('E35').Value = '  +0.11%  '
# Row 36
$ws.Range('E36').Value = '  +0.06%  '
# Row 37
$ws.Range('D37').Value = "'" + '0.587'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.82%  '
# Row 38
$ws.Range('E38').Value = '  +0.68%  '
# Row 39
$ws.Range('D39').Value = "'" + '0.825'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.72%  '
# Row 40
$ws.Range('D40').Value = "'" + '5.75'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.78%  '
# Row 41
$ws.Range('E41').Value = '  -0.09%  '
# Row 42
$ws.Range('E42').Value = '  -5.99%  '
# Row 43
$ws.Range('D43').Value = "'" + '2.15'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.65%  '
# Row 44
$ws.Range('E44').Value = '  +0.68%  '
# Row 45
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = "'" + '1.727.01'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.32%  '
# Row 46
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = "'" + '61.48'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.30%  '
# Row 47
$ws.Range('D47').Value = "'" + '86.16'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.22%  '
# Row 48
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = "'" + '0.0₆0103'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.96%  '
# Row 49
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = "'" + '1.48'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.91%  '
# Row 50
$ws.Range('D50').Value = "'" + '0.0979'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.33%  '
# Row 51
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = "'" + '0.0502'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.59%  '
